$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): shift CHARTER SCH/ELEMENTARY/HIGH SCHOOL/MIDDLE SCHL
#     one column to the right, A1:D1 -> B1:E1, carrying their style along.
$ws.Range("A1:D1").Copy()
$ws.Range("B1").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A1:D1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Clear()               # old A1 is now empty / no cell

# --- Row 2 (data): shift the four incident totals one column to the
#     right as well, A2:D2 -> B2:E2, carrying their (default) style along.
$ws.Range("A2:D2").Copy()
$ws.Range("B2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A2:D2").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats

# --- New row label in A2, styled like the row-1 header cells.
$ws.Range("A2").Value = "# of Discipline Incidents"
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
